$wb = $excel.ActiveWorkbook

# --- Notes sheet: fix "Units of measure" text ---
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("A3").Value = "Units of measure: constant 2015 US$"

# --- Data sheet: populate country-level ODA values for 2015 ---
$data = $wb.Worksheets.Item("Data")

$rows = @(
    @("AR", "Argentina", 2015, 738200),
    @("CN", "China", 2015, 37484300),
    @("IN", "India", 2015, 15930),
    @("ID", "Indonesia", 2015, 200000),
    @("JO", "Jordan", 2015, 160000),
    @("PH", "Philippines", 2015, 100000),
    @("TH", "Thailand", 2015, 4299330),
    @("VN", "Viet Nam", 2015, 2300000)
)

$r = 2
foreach ($row in $rows) {
    $data.Cells.Item($r, 1).Value = $row[0]
    $data.Cells.Item($r, 2).Value = $row[1]
    $data.Cells.Item($r, 3).Value = $row[2]
    $data.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
